$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"
$ws.Range("C2").Value = 0.016

$ws.Range("B3").Value = "''Bacteroides_coprocola_M16_DSM_17136.mat'"
$ws.Range("C3").Value = 0.003

$ws.Range("B4").Value = "''Bacteroides_coprophilus_DSM_18228.mat'"
$ws.Range("C4").Value = 0.047

$ws.Range("B5").Value = "''Bacteroides_fluxus_YIT_12057.mat'"
$ws.Range("C5").Value = 0.065

$ws.Range("B6").Value = "''Bacteroides_oleiciplenus_YIT_12058.mat'"
$ws.Range("C6").Value = 0.013

$ws.Range("B7").Value = "''Bacteroides_ovatus_ATCC_8483.mat'"
$ws.Range("C7").Value = 0.23

$ws.Range("B8").Value = "''Bacteroides_plebeius_M12_DSM_17135.mat'"
$ws.Range("C8").Value = 0

$ws.Range("B9").Value = "''Bacteroides_salyersiae_WAL_10018.mat'"
$ws.Range("C9").Value = 0.305

$ws.Range("B10").Value = "''Bacteroides_stercoris_ATCC_43183.mat'"
$ws.Range("C10").Value = 0

$ws.Range("B11").Value = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"
$ws.Range("C11").Value = -0

$ws.Range("B12").Value = "''Bacteroides_uniformis_ATCC_8492.mat'"
$ws.Range("C12").Value = 0

$ws.Range("B13").Value = "''Bacteroides_vulgatus_ATCC_8482.mat'"
$ws.Range("C13").Value = 0.322

$ws.Range("B14").Value = "''Bifidobacterium_animalis_lactis_AD011.mat'"
$ws.Range("C14").Value = 0

$ws.Range("B15").Value = "''Enterococcus_faecalis_OG1RF_ATCC_47077.mat'"
$ws.Range("C15").Value = 0

$ws.Range("B16").Value = "''Flavonifractor_plautii_ATCC_29863.mat'"
$ws.Range("C16").Value = 0

$ws.Range("B17").Value = "''Lactobacillus_plantarum_JDM1.mat'"
$ws.Range("C17").Value = 0

$ws.Range("B18").Value = "''Odoribacter_laneus_YIT_12061.mat'"
$ws.Range("C18").Value = -0

$ws.Range("B19").Value = "''Parabacteroides_distasonis_ATCC_8503.mat'"
$ws.Range("C19").Value = 0

$ws.Range("B20").Value = "''Parabacteroides_johnsonii_DSM_18315.mat'"
$ws.Range("C20").Value = 0
